$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 31441.95
$ws.Range("I21").Value = 10912.637
$ws.Range("J21").Value = 56533.332
$ws.Range("K21").Value = 10912.637
$ws.Range("L21").Value = 56533.332
$ws.Range("M21").Value = -10444.637
$ws.Range("N21").Value = -57469.332

$ws.Range("H23").Value = 31441.95
$ws.Range("I23").Value = 10912.637
$ws.Range("J23").Value = 56533.332
$ws.Range("K23").Value = 10912.637
$ws.Range("L23").Value = 56533.332
$ws.Range("M23").Value = -10678.637
$ws.Range("N23").Value = -57001.332

$ws.Range("H33").Value = 639.3684
$ws.Range("I33").Value = 547.4167
$ws.Range("J33").Value = 797
$ws.Range("K33").Value = 547.4167
$ws.Range("L33").Value = 797
$ws.Range("M33").Value = -318.4167
$ws.Range("N33").Value = -1255

$ws.Range("H58").Value = 1196.0667
$ws.Range("J58").Value = 1840.7778
$ws.Range("L58").Value = 5522.3334
$ws.Range("N58").Value = -5822.3334

$ws.Range("H70").Value = 6877.421
$ws.Range("I70").Value = 4092.3333
$ws.Range("K70").Value = 12276.9999
$ws.Range("M70").Value = -12006.9999

$ws.Range("H73").Value = 6877.421
$ws.Range("I73").Value = 4092.3333
$ws.Range("K73").Value = 12276.9999
$ws.Range("M73").Value = -11340.9999

$ws.Range("H80").Value = 5899.9375
$ws.Range("I80").Value = 1967
$ws.Range("J80").Value = 8259.700000000001
$ws.Range("K80").Value = 5901
$ws.Range("L80").Value = 24779.1
$ws.Range("M80").Value = -4903
$ws.Range("N80").Value = -26775.1

$ws.Range("H83").Value = 5899.9375
$ws.Range("I83").Value = 1967
$ws.Range("J83").Value = 8259.700000000001
$ws.Range("K83").Value = 17703
$ws.Range("L83").Value = 74337.3
$ws.Range("M83").Value = -12711
$ws.Range("N83").Value = -84321.3

$ws.Range("H100").Value = 621
$ws.Range("I100").Value = 413.8
$ws.Range("K100").Value = 413.8
$ws.Range("M100").Value = 127.2

$ws.Range("H103").Value = 832.8333
$ws.Range("I103").Value = 899.4
$ws.Range("K103").Value = 2698.2
$ws.Range("M103").Value = -2112.2

$ws.Range("H132").Value = 20410676
$ws.Range("I132").Value = 20835886
$ws.Range("K132").Value = 62507658
$ws.Range("M132").Value = -62505128

$ws.Range("H137").Value = 121491.266
$ws.Range("I137").Value = 200871.11
$ws.Range("K137").Value = 602613.33
$ws.Range("M137").Value = -600063.33

$ws.Range("H141").Value = 8522.933999999999
$ws.Range("I141").Value = 8522.933999999999
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 25568.802
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -20388.802
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2021194.9
$ws.Range("I2").Value = 2828969
$ws.Range("K2").Value = 2828969
$ws.Range("M2").Value = -2828856

$ws.Range("H32").Value = 11938.17
$ws.Range("I32").Value = 8338.344999999999
$ws.Range("J32").Value = 16287.958
$ws.Range("K32").Value = 8338.344999999999
$ws.Range("L32").Value = 16287.958
$ws.Range("M32").Value = -8051.344999999999
$ws.Range("N32").Value = -16861.958

$ws.Range("H45").Value = 4498626.5
$ws.Range("I45").Value = 5995649.5
$ws.Range("K45").Value = 5995649.5
$ws.Range("M45").Value = -5995272.5

$ws.Range("H61").Value = 3119.0557
$ws.Range("I61").Value = 2430
$ws.Range("K61").Value = 2430
$ws.Range("M61").Value = -2218

$ws.Range("H74").Value = 55286.35
$ws.Range("I74").Value = 33497.277
$ws.Range("K74").Value = 33497.277
$ws.Range("M74").Value = -32623.277

$ws.Range("H77").Value = 55286.35
$ws.Range("I77").Value = 33497.277
$ws.Range("K77").Value = 167486.385
$ws.Range("M77").Value = -163118.385

$ws.Range("H101").Value = 55000
$ws.Range("J101").Value = 55000
$ws.Range("L101").Value = 55000
$ws.Range("N101").Value = -61490

$ws.Range("H102").Value = 8337933.5
$ws.Range("I102").Value = 11907048
$ws.Range("K102").Value = 11907048
$ws.Range("M102").Value = -11905426

$ws.Range("H116").Value = 2021194.9
$ws.Range("I116").Value = 2828969
$ws.Range("K116").Value = 2828969
$ws.Range("M116").Value = -2826675

$ws.Range("H122").Value = 6571377.5
$ws.Range("I122").Value = 11697283
$ws.Range("K122").Value = 35091849
$ws.Range("M122").Value = -35089399

$ws.Range("H132").Value = 2275.524
$ws.Range("I132").Value = 2071.2
$ws.Range("K132").Value = 6213.599999999999
$ws.Range("M132").Value = -3683.599999999999

$ws.Range("H136").Value = 3119.0557
$ws.Range("I136").Value = 2430
$ws.Range("K136").Value = 7290
$ws.Range("M136").Value = -4740

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2021194.9
$ws.Range("I3").Value = 2828969
$ws.Range("K3").Value = 2828969
$ws.Range("M3").Value = -2828855

$ws.Range("H86").Value = 6683004
$ws.Range("I86").Value = 8344755
$ws.Range("K86").Value = 8344755
$ws.Range("M86").Value = -8343632

$ws.Range("H89").Value = 6683004
$ws.Range("I89").Value = 8344755
$ws.Range("K89").Value = 41723775
$ws.Range("M89").Value = -41718159

$ws.Range("H134").Value = 5584.222
$ws.Range("I134").Value = 1964.8
$ws.Range("K134").Value = 5894.4
$ws.Range("M134").Value = -3359.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 40150.043
$ws.Range("J31").Value = 152898
$ws.Range("L31").Value = 152898
$ws.Range("N31").Value = -153488

$ws.Range("H34").Value = 40150.043
$ws.Range("J34").Value = 152898
$ws.Range("L34").Value = 152898
$ws.Range("N34").Value = -153302

$ws.Range("H107").Value = 32260400
$ws.Range("I107").Value = 1850.6842
$ws.Range("J107").Value = 83336430
$ws.Range("K107").Value = 1850.6842
$ws.Range("L107").Value = 83336430
$ws.Range("M107").Value = 69.31580000000008
$ws.Range("N107").Value = -83340270

$ws.Range("H132").Value = 57147.207
$ws.Range("I132").Value = 36104.035
$ws.Range("J132").Value = 179197.6
$ws.Range("K132").Value = 108312.105
$ws.Range("L132").Value = 537592.8
$ws.Range("M132").Value = -105782.105
$ws.Range("N132").Value = -542652.8

$ws.Range("H134").Value = 22709.408
$ws.Range("I134").Value = 30897.033
$ws.Range("K134").Value = 92691.099
$ws.Range("M134").Value = -90156.099

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 7839.643
$ws.Range("I11").Value = 5362
$ws.Range("K11").Value = 16086
$ws.Range("M11").Value = -15946

$ws.Range("H119").Value = 5841.4
$ws.Range("J119").Value = 3766
$ws.Range("L119").Value = 11298
$ws.Range("N119").Value = -20974

$ws.Range("H129").Value = 1053391
$ws.Range("I129").Value = 1250712.2
$ws.Range("J129").Value = 1011
$ws.Range("K129").Value = 3752136.6
$ws.Range("L129").Value = 3033
$ws.Range("M129").Value = -3747136.6
$ws.Range("N129").Value = -13033

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H132").Value = 3058.3235
$ws.Range("I132").Value = 2862.5173
$ws.Range("K132").Value = 8587.5519
$ws.Range("M132").Value = -6057.5519

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 1000
$ws.Range("I18").Value = 1000
$ws.Range("K18").Value = 1000
$ws.Range("M18").Value = -828

$ws.Range("H40").Value = 8899.6
$ws.Range("I40").Value = 4833.3335
$ws.Range("K40").Value = 4833.3335
$ws.Range("M40").Value = -4697.3335

$ws.Range("H51").Value = 42250
$ws.Range("J51").Value = 42250
$ws.Range("L51").Value = 42250
$ws.Range("N51").Value = -43206

$ws.Range("H55").Value = 1384.5238
$ws.Range("I55").Value = 1317.6154
$ws.Range("K55").Value = 1317.6154
$ws.Range("M55").Value = -1144.6154

$ws.Range("H100").Value = 3540
$ws.Range("I100").Value = 3301.5
$ws.Range("K100").Value = 3301.5
$ws.Range("M100").Value = -2760.5

$ws.Range("H119").Value = 105000
$ws.Range("J119").Value = 105000
$ws.Range("L119").Value = 105000
$ws.Range("N119").Value = -114676

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 38566
$ws.Range("J95").Value = 38566
$ws.Range("L95").Value = 38566
$ws.Range("N95").Value = -44058
